$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").Value = ""
$ws.Range("H33").Value = 95.07692
$ws.Range("I33").Value = 53.88889
$ws.Range("K33").Value = 53.88889
$ws.Range("M33").Value = 175.11111
$ws.Range("H132").Value = 1277.5714
$ws.Range("I132").Value = 1277.5714
$ws.Range("K132").Value = 3832.7142
$ws.Range("M132").Value = -1302.7142
$ws.Range("H138").Value = 4111.362
$ws.Range("J138").Value = 4364.381
$ws.Range("L138").Value = 13093.143
$ws.Range("N138").Value = -23373.143
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1999.7142
$ws.Range("I2").Value = 2000
$ws.Range("K2").Value = 2000
$ws.Range("M2").Value = -1887
$ws.Range("H32").Value = 9099.706
$ws.Range("I32").Value = 6979.6665
$ws.Range("K32").Value = 6979.6665
$ws.Range("M32").Value = -6692.6665
$ws.Range("H33").Value = 1800
$ws.Range("I33").Value = 1800
$ws.Range("K33").Value = 1800
$ws.Range("M33").Value = -1471
$ws.Range("H45").Value = 1749.4445
$ws.Range("I45").Value = 1749.4445
$ws.Range("K45").Value = 1749.4445
$ws.Range("M45").Value = -1372.4445
$ws.Range("H61").Value = 2133.3333
$ws.Range("I61").Value = 2133.3333
$ws.Range("K61").Value = 2133.3333
$ws.Range("M61").Value = -1921.3333
$ws.Range("H116").Value = 1999.7142
$ws.Range("I116").Value = 2000
$ws.Range("K116").Value = 2000
$ws.Range("M116").Value = 294
$ws.Range("H136").Value = 2133.3333
$ws.Range("I136").Value = 2133.3333
$ws.Range("K136").Value = 6399.999899999999
$ws.Range("M136").Value = -3849.999899999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1999.7142
$ws.Range("I3").Value = 2000
$ws.Range("K3").Value = 2000
$ws.Range("M3").Value = -1886
$ws.Range("H99").Value = 3998.4
$ws.Range("I99").Value = 3833.0667
$ws.Range("K99").Value = 3833.0667
$ws.Range("M99").Value = -2335.0667
$ws.Range("H105").Value = 3926.739
$ws.Range("I105").Value = 3024.1052
$ws.Range("K105").Value = 3024.1052
$ws.Range("M105").Value = -1277.1052
$ws.Range("H134").Value = 2115.3635
$ws.Range("I134").Value = 2151.9
$ws.Range("K134").Value = 6455.700000000001
$ws.Range("M134").Value = -3920.700000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 46943.555
$ws.Range("I62").Value = 2624.75
$ws.Range("J62").Value = 82398.60000000001
$ws.Range("K62").Value = 2624.75
$ws.Range("L62").Value = 82398.60000000001
$ws.Range("M62").Value = -2000.75
$ws.Range("N62").Value = -83646.60000000001
$ws.Range("H65").Value = 46943.555
$ws.Range("I65").Value = 2624.75
$ws.Range("J65").Value = 82398.60000000001
$ws.Range("K65").Value = 13123.75
$ws.Range("L65").Value = 411993
$ws.Range("M65").Value = -10003.75
$ws.Range("N65").Value = -418233
$ws.Range("H99").Value = 16511.773
$ws.Range("I99").Value = 14078.7
$ws.Range("J99").Value = 18539.334
$ws.Range("K99").Value = 14078.7
$ws.Range("L99").Value = 18539.334
$ws.Range("M99").Value = -12580.7
$ws.Range("N99").Value = -21535.334
$ws.Range("H107").Value = 584.96295
$ws.Range("I107").Value = 425.57895
$ws.Range("K107").Value = 425.57895
$ws.Range("M107").Value = 1494.42105
$ws.Range("H126").Value = 16511.773
$ws.Range("I126").Value = 14078.7
$ws.Range("J126").Value = 18539.334
$ws.Range("K126").Value = 42236.10000000001
$ws.Range("L126").Value = 55618.00199999999
$ws.Range("M126").Value = -39766.10000000001
$ws.Range("N126").Value = -60558.00199999999
$ws.Range("H132").Value = 2828.9333
$ws.Range("J132").Value = 4072.4285
$ws.Range("L132").Value = 12217.2855
$ws.Range("N132").Value = -17277.2855
$ws.Range("H141").Value = 15499.5
$ws.Range("J141").Value = 15499.5
$ws.Range("L141").Value = 15499.5
$ws.Range("N141").Value = -25859.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 5000
$ws.Range("I31").Value = 5000
$ws.Range("K31").Value = 5000
$ws.Range("M31").Value = -4708
$ws.Range("H35").Value = 5000000
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").Value = ""
$ws.Range("H37").Value = 5000
$ws.Range("I37").Value = 5000
$ws.Range("K37").Value = 5000
$ws.Range("M37").Value = -4723
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = ""
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").Value = ""
$ws.Range("H97").Value = 854.64703
$ws.Range("I97").Value = 730.38464
$ws.Range("J97").Value = 1258.5
$ws.Range("K97").Value = 730.38464
$ws.Range("L97").Value = 1258.5
$ws.Range("M97").Value = -234.38464
$ws.Range("N97").Value = -2250.5
$ws.Range("H122").Value = 58628.277
$ws.Range("I122").Value = 2554.3333
$ws.Range("J122").Value = 338998
$ws.Range("K122").Value = 7662.999899999999
$ws.Range("L122").Value = 1016994
$ws.Range("M122").Value = -5212.999899999999
$ws.Range("N122").Value = -1021894
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2950
$ws.Range("J46").Value = 4400
$ws.Range("L46").Value = 4400
$ws.Range("N46").Value = -4776
$ws.Range("H93").Value = 1000
$ws.Range("I93").Value = 1000
$ws.Range("K93").Value = 1000
$ws.Range("M93").Value = 248
$ws.Range("H122").Value = 9334.333000000001
$ws.Range("I122").Value = 9334.333000000001
$ws.Range("K122").Value = 28002.999
$ws.Range("M122").Value = -25552.999
$ws.Range("H136").Value = 2491.25
$ws.Range("I136").Value = 1957.5
$ws.Range("J136").Value = 3025
$ws.Range("K136").Value = 5872.5
$ws.Range("L136").Value = 9075
$ws.Range("M136").Value = -3322.5
$ws.Range("N136").Value = -14175
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 15000
$ws.Range("J81").Value = 16333.333
$ws.Range("L81").Value = 32666.666
$ws.Range("N81").Value = -34788.666
$ws.Range("H84").Value = 15000
$ws.Range("J84").Value = 16333.333
$ws.Range("L84").Value = 163333.33
$ws.Range("N84").Value = -173941.33
$ws.Range("H107").Value = 642.9286
$ws.Range("J107").Value = 757.44446
$ws.Range("L107").Value = 2272.33338
$ws.Range("N107").Value = -6112.33338
$ws.Range("H132").Value = 2727
$ws.Range("I132").Value = 2340.75
$ws.Range("K132").Value = 7022.25
$ws.Range("M132").Value = -4492.25
